# Updated cryptos list on Thu Aug 22 15:48:20 UTC 2024 with GitHub Actions
#
# Applies the latest scraped Price / Volume(1h) figures to the "cryptos"
# sheet, and re-sorts two adjacent coin pairs whose relative ranking
# changed (PEPE <-> PancakeSwap around rows 29/30, and RenderToken <->
# EnergySwap around rows 46/47).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Bitcoin ---------------------------------------------------
$ws.Range("D2").Value = "60.095.52"
$ws.Range("E2").Value = "  +0.89%  "

# --- Row 3: Ethereum ---------------------------------------------------
$ws.Range("D3").Value = "2.592.67"
$ws.Range("E3").Value = "  +0.08%  "

# --- Row 4: TetherUSD ---------------------------------------------------
$ws.Range("E4").Value = "  +0.06%  "

# --- Row 5: BNB ---------------------------------------------------
$ws.Range("D5").Value = "577.99"
$ws.Range("E5").Value = "  +4.21%  "

# --- Row 6: Solana ---------------------------------------------------
$ws.Range("D6").Value = "142.41"
$ws.Range("E6").Value = "  +1.22%  "

# --- Row 8: XRP ---------------------------------------------------
$ws.Range("E8").Value = "  +0.35%  "

# --- Row 9: LidoStakedEther ---------------------------------------------------
$ws.Range("D9").Value = "2.599.21"
$ws.Range("E9").Value = "  -0.28%  "

# --- Row 10: Toncoin ---------------------------------------------------
$ws.Range("E10").Value = "  -2.86%  "

# --- Row 11: Dogecoin ---------------------------------------------------
$ws.Range("E11").Value = "  +0.82%  "

# --- Row 12: TRON ---------------------------------------------------
$ws.Range("E12").Value = "  -2.77%  "

# --- Row 13: Cardano ---------------------------------------------------
$ws.Range("E13").Value = "  +3.43%  "

# --- Row 14: WrappedliquidstakedEther2.0 ---------------------------------------------------
$ws.Range("D14").Value = "3.057.72"
$ws.Range("E14").Value = "  +0.24%  "

# --- Row 15: Avalanche ---------------------------------------------------
$ws.Range("D15").Value = "24.63"
$ws.Range("E15").Value = "  +6.83%  "

# --- Row 16: WrappedBTC ---------------------------------------------------
$ws.Range("D16").Value = "60.104.35"
$ws.Range("E16").Value = "  +0.97%  "

# --- Row 17: ShibaInu ---------------------------------------------------
$ws.Range("D17").Value = "0.0000141"
$ws.Range("E17").Value = "  +2.42%  "

# --- Row 18: WrappedEther ---------------------------------------------------
$ws.Range("D18").Value = "2.601.17"
$ws.Range("E18").Value = "  +0.03%  "

# --- Row 19: Chainlink ---------------------------------------------------
$ws.Range("D19").Value = "11.48"
$ws.Range("E19").Value = "  +10.03%  "

# --- Row 20: Polkadot ---------------------------------------------------
$ws.Range("D20").Value = "4.63"
$ws.Range("E20").Value = "  +1.53%  "

# --- Row 21: BitcoinCash ---------------------------------------------------
$ws.Range("D21").Value = "345.96"
$ws.Range("E21").Value = "  +1.69%  "

# --- Row 22: Uniswap ---------------------------------------------------
$ws.Range("D22").Value = "6.88"
$ws.Range("E22").Value = "  +4.50%  "

# --- Row 23: Dai ---------------------------------------------------
$ws.Range("E23").Value = "  +0.28%  "

# --- Row 24: Polygon ---------------------------------------------------
$ws.Range("D24").Value = "0.525"
$ws.Range("E24").Value = "  +8.21%  "

# --- Row 25: Litecoin ---------------------------------------------------
$ws.Range("D25").Value = "62.96"
$ws.Range("E25").Value = "  -0.10%  "

# --- Row 26: Binance-PegBSC-USD ---------------------------------------------------
$ws.Range("E26").Value = "  +0.31%  "

# --- Row 27: Kaspa ---------------------------------------------------
$ws.Range("E27").Value = "  +0.09%  "

# --- Row 28: InternetComputer(DFINITY) ---------------------------------------------------
$ws.Range("D28").Value = "8.02"
$ws.Range("E28").Value = "  +7.11%  "

# --- Rows 29/30: PEPE and PancakeSwap swap ranking order ----------------
# Row 29 becomes PancakeSwap (was PEPE), Row 30 becomes PEPE (was PancakeSwap)
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "1.87"
$ws.Range("E29").Value = "  +10.87%  "

$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0" + [char]0x2083 + "0788"
$ws.Range("E30").Value = "  +2.29%  "

# --- Row 31: USDe ---------------------------------------------------
$ws.Range("D31").Value = "0.998"

# --- Row 32: Aptos ---------------------------------------------------
$ws.Range("E32").Value = "  +3.83%  "

# --- Row 33: Monero ---------------------------------------------------
$ws.Range("D33").Value = "164.27"
$ws.Range("E33").Value = "  +4.30%  "

# --- Row 34: EthereumClassic ---------------------------------------------------
$ws.Range("E34").Value = "  +0.03%  "

# --- Row 35: NEARProtocol ---------------------------------------------------
$ws.Range("D35").Value = "4.27"
$ws.Range("E35").Value = "  +4.86%  "

# --- Row 36: Fetch.AI ---------------------------------------------------
$ws.Range("D36").Value = "0.985"
$ws.Range("E36").Value = "  +6.99%  "

# --- Row 37: ImmutableX ---------------------------------------------------
$ws.Range("E37").Value = "  +6.72%  "

# --- Row 38: Stacks ---------------------------------------------------
$ws.Range("E38").Value = "  +8.86%  "

# --- Row 39: OKB ---------------------------------------------------
$ws.Range("D39").Value = "37.96"
$ws.Range("E39").Value = "  +1.04%  "

# --- Row 40: Filecoin ---------------------------------------------------
$ws.Range("D40").Value = "3.89"
$ws.Range("E40").Value = "  +5.85%  "

# --- Row 41: Bittensor ---------------------------------------------------
$ws.Range("D41").Value = "307.85"
$ws.Range("E41").Value = "  +5.98%  "

# --- Row 42: SuiNetwork ---------------------------------------------------
$ws.Range("E42").Value = "  -0.52%  "

# --- Row 43: Aave ---------------------------------------------------
$ws.Range("D43").Value = "135.26"
$ws.Range("E43").Value = "  -1.04%  "

# --- Row 44: FirstDigitalUSD ---------------------------------------------------
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.31%  "

# --- Row 45: Stellar ---------------------------------------------------
$ws.Range("D45").Value = "0.0987"
$ws.Range("E45").Value = "  +1.33%  "

# --- Rows 46/47: RenderToken and EnergySwap swap ranking order ----------
# Row 46 becomes EnergySwap (was RenderToken), Row 47 becomes RenderToken (was EnergySwap)
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "19.68"
$ws.Range("E46").Value = "  +3.75%  "

$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "4.99"
$ws.Range("E47").Value = "  +9.56%  "

# --- Row 48: Mantle ---------------------------------------------------
$ws.Range("E48").Value = "  +0.50%  "

# --- Row 49: Hedera ---------------------------------------------------
$ws.Range("D49").Value = "0.0547"
$ws.Range("E49").Value = "  +2.10%  "

# --- Row 50: InjectiveProtocol ---------------------------------------------------
$ws.Range("D50").Value = "19.98"
$ws.Range("E50").Value = "  +6.90%  "

# --- Row 51: VeChain ---------------------------------------------------
$ws.Range("E51").Value = "  +2.16%  "
